$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table has 20 rows, but only every 4th row (1, 5, 9, 13, 17) holds the
# 5 division-problem cells; the rows in between are blank spacer rows.
# Addressing cells by (row, column) avoids any ambiguity from Find/Replace
# when new values momentarily collide with not-yet-processed old values.
$contentRows = @(1, 5, 9, 13, 17)

$newValues = @(
    @("85÷5=", "24÷3=", "38÷2=", "61÷3=", "94÷3="),
    @("12÷6=", "98÷7=", "47÷3=", "32÷7=", "81÷4="),
    @("36÷4=", "42÷8=", "83÷3=", "66÷7=", "10÷4="),
    @("50÷7=", "90÷6=", "67÷7=", "26÷9=", "24÷8="),
    @("68÷8=", "42÷7=", "19÷3=", "92÷7=", "91÷4=")
)

for ($i = 0; $i -lt $contentRows.Length; $i++) {
    $row = $contentRows[$i]
    $values = $newValues[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
